$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 230, shifting the existing rows 230-232 down to 231-233.
$ws.Rows.Item(230).Insert()

# Populate the newly inserted row 230 with the new weekly record.
$ws.Range("A230").Value2 = 5
$ws.Range("B230").Value2 = "Macroferia Regional de Talca"
$ws.Range("C230").Value2 = "Maule"
$ws.Range("D230").Value2 = 44656
$ws.Range("E230").Value2 = 7
$ws.Range("F230").Value2 = 100112021
$ws.Range("G230").Value2 = "Ají"
$ws.Range("H230").Value2 = "Cristal"
$ws.Range("I230").Value2 = "Primera"
$ws.Range("J230").Value2 = 200
$ws.Range("K230").Value2 = 14000
$ws.Range("L230").Value2 = 14000
$ws.Range("M230").Value2 = 14000
$ws.Range("N230").Value2 = "$/saco 25 kilos"
$ws.Range("O230").Value2 = "Región del Maule"
$ws.Range("P230").Value2 = 560
$ws.Range("Q230").Value2 = 25
$ws.Range("R230").Value2 = "Hortaliza"
